$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2:A5 with the consolidated token descriptions
$ws.Range("A2").Value = "('Bird Soldier', ['Token Creature — Bird Soldier', 'Flying', '1/1'])"
$ws.Range("A3").Value = "('Dragon', ['Token Creature — Dragon', 'Flying, devour 2', '1/1'])"
$ws.Range("A4").Value = "('Lizard', ['Token Creature — Lizard', '2/2'])"
$ws.Range("A5").Value = "('Zombie Wizard', ['Token Creature — Zombie Wizard', '1/1'])"

# Remove the now-unused rows 6 through 15 (their content has been folded into A2:A5)
$ws.Range("A6:A15").ClearContents()
